$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 12073.454
$ws.Range("I98").Value = 15976.25
$ws.Range("K98").Value = 15976.25
$ws.Range("M98").Value = -14478.25
$ws.Range("H107").Value = 1836.409
$ws.Range("I107").Value = 1566.7222
$ws.Range("K107").Value = 1566.7222
$ws.Range("M107").Value = 353.2778000000001
$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680
$ws.Range("H112").Value = 2713.6365
$ws.Range("J112").Value = 2713.6365
$ws.Range("L112").Value = 8140.9095
$ws.Range("N112").Value = -10356.9095
$ws.Range("H122").Value = 12073.454
$ws.Range("I122").Value = 15976.25
$ws.Range("K122").Value = 47928.75
$ws.Range("M122").Value = -45478.75
$ws.Range("H129").Value = 933.2
$ws.Range("I129").Value = 797
$ws.Range("J129").Value = 991.5714
$ws.Range("K129").Value = 2391
$ws.Range("L129").Value = 2974.7142
$ws.Range("M129").Value = 2609
$ws.Range("N129").Value = -12974.7142
$ws.Range("H132").Value = 7414457.5
$ws.Range("I132").Value = 8134407
$ws.Range("K132").Value = 24403221
$ws.Range("M132").Value = -24400691
$ws.Range("H137").Value = 1764.2273
$ws.Range("I137").Value = 1690
$ws.Range("J137").Value = 1838.4546
$ws.Range("K137").Value = 5070
$ws.Range("L137").Value = 5515.3638
$ws.Range("M137").Value = -2520
$ws.Range("N137").Value = -10615.3638
$ws.Range("H138").Value = 3092.76
$ws.Range("J138").Value = 3097.923
$ws.Range("L138").Value = 9293.769
$ws.Range("N138").Value = -19573.769

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22456.834
$ws.Range("I32").Value = 14079.257
$ws.Range("K32").Value = 14079.257
$ws.Range("M32").Value = -13792.257
$ws.Range("H61").Value = 40001492
$ws.Range("I61").Value = 43479580
$ws.Range("K61").Value = 43479580
$ws.Range("M61").Value = -43479368
$ws.Range("H74").Value = 2170.7307
$ws.Range("I74").Value = 1160.9
$ws.Range("K74").Value = 1160.9
$ws.Range("M74").Value = -286.9000000000001
$ws.Range("H77").Value = 2170.7307
$ws.Range("I77").Value = 1160.9
$ws.Range("K77").Value = 5804.5
$ws.Range("M77").Value = -1436.5
$ws.Range("H136").Value = 40001492
$ws.Range("I136").Value = 43479580
$ws.Range("K136").Value = 130438740
$ws.Range("M136").Value = -130436190

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1897.24
$ws.Range("I31").Value = 1846.8
$ws.Range("J31").Value = 2099
$ws.Range("K31").Value = 1846.8
$ws.Range("L31").Value = 2099
$ws.Range("M31").Value = -1551.8
$ws.Range("N31").Value = -2689
$ws.Range("H34").Value = 1897.24
$ws.Range("I34").Value = 1846.8
$ws.Range("J34").Value = 2099
$ws.Range("K34").Value = 1846.8
$ws.Range("L34").Value = 2099
$ws.Range("M34").Value = -1644.8
$ws.Range("N34").Value = -2503
$ws.Range("H58").Value = 3223.2449
$ws.Range("I58").Value = 817.3103599999999
$ws.Range("K58").Value = 817.3103599999999
$ws.Range("M58").Value = -614.3103599999999
$ws.Range("H105").Value = 810.9231
$ws.Range("I105").Value = 754.7
$ws.Range("J105").Value = 998.3333
$ws.Range("K105").Value = 754.7
$ws.Range("L105").Value = 998.3333
$ws.Range("M105").Value = 992.3
$ws.Range("N105").Value = -4492.3333
$ws.Range("H122").Value = 1999.5
$ws.Range("I122").Value = 1999
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5997
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3547
$ws.Range("N122").Value = -10900
$ws.Range("H136").Value = 3223.2449
$ws.Range("I136").Value = 817.3103599999999
$ws.Range("K136").Value = 2451.93108
$ws.Range("M136").Value = 98.06892000000016
$ws.Range("H141").Value = 486330.47
$ws.Range("J141").Value = 486330.47
$ws.Range("L141").Value = 486330.47
$ws.Range("N141").Value = -496690.47

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 506
$ws.Range("I13").Value = 335.8
$ws.Range("J13").Value = 789.6667
$ws.Range("K13").Value = 1007.4
$ws.Range("L13").Value = 2369.0001
$ws.Range("M13").Value = -839.4000000000001
$ws.Range("N13").Value = -2705.0001
$ws.Range("H34").Value = 7693563
$ws.Range("I34").Value = 188.85715
$ws.Range("J34").Value = 16669167
$ws.Range("K34").Value = 566.5714499999999
$ws.Range("L34").Value = 50007501
$ws.Range("M34").Value = -482.5714499999999
$ws.Range("N34").Value = -50007669
$ws.Range("H39").Value = 4254.5
$ws.Range("J39").Value = 4345.4
$ws.Range("L39").Value = 13036.2
$ws.Range("N39").Value = -13624.2
$ws.Range("H55").Value = 2876
$ws.Range("J55").Value = 3500
$ws.Range("L55").Value = 10500
$ws.Range("N55").Value = -10854
$ws.Range("H107").Value = 5675.68
$ws.Range("J107").Value = 9230.200000000001
$ws.Range("L107").Value = 27690.6
$ws.Range("N107").Value = -31530.6
$ws.Range("H113").Value = 716.46875
$ws.Range("I113").Value = 600
$ws.Range("K113").Value = 1800
$ws.Range("M113").Value = 370
$ws.Range("H126").Value = 4377.8
$ws.Range("I126").Value = 2963.3333
$ws.Range("K126").Value = 8889.999899999999
$ws.Range("M126").Value = -3949.999899999999
$ws.Range("H131").Value = 18897200
$ws.Range("J131").Value = 33663.023
$ws.Range("L131").Value = 100989.069
$ws.Range("N131").Value = -111069.069
$ws.Range("H140").Value = 22983.244
$ws.Range("I140").Value = 47696.684
$ws.Range("K140").Value = 143090.052
$ws.Range("M140").Value = -137910.052

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3880.1667
$ws.Range("I70").Value = 3815.4285
$ws.Range("J70").Value = 3970.8
$ws.Range("K70").Value = 3815.4285
$ws.Range("L70").Value = 3970.8
$ws.Range("M70").Value = -3545.4285
$ws.Range("N70").Value = -4510.8
$ws.Range("H73").Value = 3880.1667
$ws.Range("I73").Value = 3815.4285
$ws.Range("J73").Value = 3970.8
$ws.Range("K73").Value = 3815.4285
$ws.Range("L73").Value = 3970.8
$ws.Range("M73").Value = -2879.4285
$ws.Range("N73").Value = -5842.8
$ws.Range("H132").Value = 7984.8096
$ws.Range("I132").Value = 9609.4
$ws.Range("J132").Value = 3923.3333
$ws.Range("K132").Value = 28828.2
$ws.Range("L132").Value = 11769.9999
$ws.Range("M132").Value = -26298.2
$ws.Range("N132").Value = -16829.9999

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3007.5
$ws.Range("I40").Value = 3012.2222
$ws.Range("K40").Value = 3012.2222
$ws.Range("M40").Value = -2876.2222
$ws.Range("H122").Value = 94445110
$ws.Range("I122").Value = 141666670
$ws.Range("K122").Value = 425000010
$ws.Range("M122").Value = -424997560
$ws.Range("H132").Value = 5356.7144
$ws.Range("I132").Value = 5999.5
$ws.Range("K132").Value = 17998.5
$ws.Range("M132").Value = -15468.5
$ws.Range("H136").Value = 2781.6667
$ws.Range("I136").Value = 2672.5
$ws.Range("K136").Value = 8017.5
$ws.Range("M136").Value = -5467.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3981.5
$ws.Range("I132").Value = 3671
$ws.Range("K132").Value = 11013
$ws.Range("M132").Value = -8483
$ws.Range("H136").Value = 2124.6875
$ws.Range("I136").Value = 2099.0908
$ws.Range("K136").Value = 6297.2724
$ws.Range("M136").Value = -3747.2724
